$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sale price amount" / "Sale price currency" columns (P:Q) used to live
# on this sales import template; they now belong on the contract file
# template instead, so remove them here together with the now-unused
# trailing placeholder rows that followed the header table.
$ws.Rows("11:54").Delete()
$ws.Columns("P:Q").Delete()

$ws.Range("P4").Select() | Out-Null
